$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column C (header through total) into column D to replicate formatting exactly
$ws.Range("C4:C15").Copy()
$ws.Range("D4:D15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new header text
$ws.Range("D4").Value = "#3"

# New data value for D5, rest stay blank
$ws.Range("D5").Value = 0.5

# Total row D15 - sum formula
$ws.Range("D15").Formula = "=SUM(D5:D14)"

# Update selection to match the target view
$ws.Range("D6:D14").Select()
